$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.963.96"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "1.919.91"
$ws.Range("E3").Value = "  +1.26%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "320.53"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "0.5057"
$ws.Range("E7").Value = "  -2.03%  "

$ws.Range("D8").Value = "0.4044"
$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("D9").Value = "0.08298"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").Value = "1.104"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").Value = "41.94"
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").Value = "24.01"
$ws.Range("E12").Value = "  +2.94%  "

$ws.Range("D13").Value = "1.919.78"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").Value = "6.403"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "'7.240"
$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").Value = "0.06502"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").Value = "18.19"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "'5.950"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "29.993.50"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").Value = "2.196"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("D27").Value = "2.140.48"
$ws.Range("E27").Value = "  +1.33%  "

$ws.Range("D28").Value = "'162.50"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "2.312"
$ws.Range("E29").Value = "  -1.77%  "

$ws.Range("D30").Value = "129.33"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").Value = "1.137"
$ws.Range("E31").Value = "  +3.79%  "

$ws.Range("E32").Value = "  -1.91%  "

$ws.Range("D33").Value = "'5.970"
$ws.Range("E33").Value = "  -2.20%  "

$ws.Range("D34").Value = "3.837"
$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("D35").Value = "0.02461"
$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("D36").Value = "5.411"
$ws.Range("E36").Value = "  +2.45%  "

$ws.Range("D37").Value = "0.06424"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("D38").Value = "0.2155"
$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "8.733"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.196"
$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("D41").Value = "0.6465"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").Value = "11.37"
$ws.Range("E42").Value = "  -3.23%  "

$ws.Range("D43").Value = "1.214"
$ws.Range("E43").Value = "  -1.24%  "

$ws.Range("D44").Value = "2.224"
$ws.Range("E44").Value = "  +8.30%  "

$ws.Range("D45").Value = "'13.30"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("D46").Value = "0.6048"
$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("D47").Value = "3.636"
$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "122.26"
$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.208"
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("D50").Value = "'79.10"
$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").Value = "'1.130"
$ws.Range("E51").Value = "  -2.81%  "
